$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.673.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.998.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.07%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.995.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.111"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.70%  "
$ws.Range("E11").Value = "  -4.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.360"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.517.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.123"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.769.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.996.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.92%  "
$ws.Range("E18").Value = "  -4.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "388.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.459"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.49%  "
$ws.Range("E26").Value = "  -8.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0958"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.04%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.513.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.75%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.13%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.658"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0593"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0246"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.25%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "270.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0937"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.92%  "
